$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.561
$ws.Range("B12").Value = 5.325
$ws.Range("D14").Value = -7.547
$ws.Range("D26").Value = -8.119
$ws.Range("D31").Value = -8.414000000000001
$ws.Range("B32").Value = 6.406999999999999
$ws.Range("D35").Value = -7.672
$ws.Range("B36").Value = 8.548
$ws.Range("D37").Value = -7.741
$ws.Range("B38").Value = 5.438000000000001
$ws.Range("D45").Value = -7.498
$ws.Range("B46").Value = 6.377000000000001
$ws.Range("B54").Value = 5.154999999999999
$ws.Range("B55").Value = 4.684
$ws.Range("D57").Value = -8.101000000000001
$ws.Range("B67").Value = 5.286
$ws.Range("B69").Value = 5.135999999999999
$ws.Range("B72").Value = 5.380999999999999
$ws.Range("B91").Value = 6.378
$ws.Range("B99").Value = 5.217
$ws.Range("D100").Value = -8.280000000000001
$ws.Range("D102").Value = -7.861
